$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 12 (11th project entry) with the new YoloV3 object detection entry
$ws.Range("C12").Value = "Masked Language Modeling using BERT mutiple model"
$ws.Range("D12").Value = "BERT"
$ws.Range("E12").Value = "Done"

# Copy the "Done" status cell style (green fill) from an existing row (E11) into E12
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to C16, matching the author's last selected cell
$ws.Range("C16").Select() | Out-Null
